# TC27_Verify_Store_room.xlsx edit — "Logic change for Logged in User"
# 1. Insert a new row 4: CLICK / LoginURL / CSS  (click the Login URL before entering creds)
#    (this naturally shifts every row below it down by one, including the final
#     Logout_RegisteredUser row which becomes row 28 automatically)
# 2. Rename Uname1/Password1/LoginButton1 -> Uname/Password/LoginButton (now at rows 5-7)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC27_Verify_Store_room")

# --- Insert new row 4 (shifts old rows 4.. down by one) ---
$ws.Rows.Item(4).Insert()

# --- Rename the login-credentials rows (now rows 5, 6, 7) first so the new
#     "LoginButton" shared string gets registered before "LoginURL" does,
#     matching the target shared-strings table order ---
$ws.Range("C5").Value = "Uname"
$ws.Range("C6").Value = "Password"
$ws.Range("C7").Value = "LoginButton"

# --- Fill in the newly inserted row 4 ---
$ws.Range("B4").Value = "CLICK"
$ws.Range("C4").Value = "LoginURL"
$ws.Range("D4").Value = "CSS"

# --- Update the sheet view: drop the frozen top-left scroll position and
#     select B4:D4 (the newly inserted row) instead of the old C27 selection ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("B4:D4").Select()
